$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 8396.5
$ws.Range("I40").Value = 5994
$ws.Range("J40").Value = 8997.125
$ws.Range("K40").Value = 5994
$ws.Range("L40").Value = 8997.125
$ws.Range("M40").Value = -5819
$ws.Range("N40").Value = -9347.125
# Row 48
$ws.Range("H48").Value = 487.5
$ws.Range("I48").Value = 487.5
$ws.Range("K48").Value = 1462.5
$ws.Range("M48").Value = -1170.5
# Row 56
$ws.Range("H56").Value = 487.5
$ws.Range("I56").Value = 487.5
$ws.Range("K56").Value = 1462.5
$ws.Range("M56").Value = -928.5
# Row 64
$ws.Range("H64").Value = 3250
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4496
# Row 67
$ws.Range("H67").Value = 3250
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5716
# Row 69
$ws.Range("H69").Value = 2833.3333
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 3250
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 9750
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -11498
# Row 72
$ws.Range("H72").Value = 2833.3333
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 3250
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 29250
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -37986
# Row 137
$ws.Range("H137").Value = 3120.9473
$ws.Range("I137").Value = 2659.6
$ws.Range("K137").Value = 7978.799999999999
$ws.Range("M137").Value = -5428.799999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 39
$ws.Range("H39").Value = 2516
$ws.Range("I39").Value = 2516
$ws.Range("K39").Value = 2516
$ws.Range("M39").Value = -1996

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 10112359
$ws.Range("I22").Value = 10112359
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10112359
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -10112186
$ws.Range("N22").ClearContents()
# Row 64
$ws.Range("H64").Value = 807.5
$ws.Range("I64").Value = 861.6667
$ws.Range("K64").Value = 861.6667
$ws.Range("M64").Value = -636.6667
# Row 67
$ws.Range("H67").Value = 807.5
$ws.Range("I67").Value = 861.6667
$ws.Range("K67").Value = 861.6667
$ws.Range("M67").Value = -81.66669999999999
# Row 80
$ws.Range("H80").Value = 678.8
$ws.Range("I80").Value = 607.5
$ws.Range("K80").Value = 607.5
$ws.Range("M80").Value = 390.5
# Row 83
$ws.Range("H83").Value = 678.8
$ws.Range("I83").Value = 607.5
$ws.Range("K83").Value = 3037.5
$ws.Range("M83").Value = 1954.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2574
# Row 31
$ws.Range("H31").Value = 6174.75
$ws.Range("I31").Value = 5771.643
$ws.Range("K31").Value = 5771.643
$ws.Range("M31").Value = -5476.643
# Row 34
$ws.Range("H34").Value = 6174.75
$ws.Range("I34").Value = 5771.643
$ws.Range("K34").Value = 5771.643
$ws.Range("M34").Value = -5569.643
# Row 62
$ws.Range("H62").Value = 8417.154
$ws.Range("I62").Value = 9544.833000000001
$ws.Range("J62").Value = 7450.5713
$ws.Range("K62").Value = 9544.833000000001
$ws.Range("L62").Value = 7450.5713
$ws.Range("M62").Value = -8920.833000000001
$ws.Range("N62").Value = -8698.5713
# Row 65
$ws.Range("H65").Value = 8417.154
$ws.Range("I65").Value = 9544.833000000001
$ws.Range("J65").Value = 7450.5713
$ws.Range("K65").Value = 47724.165
$ws.Range("L65").Value = 37252.85649999999
$ws.Range("M65").Value = -44604.165
$ws.Range("N65").Value = -43492.85649999999
# Row 99
$ws.Range("H99").Value = 1187.5
$ws.Range("J99").Value = 1375
$ws.Range("L99").Value = 1375
$ws.Range("N99").Value = -4371
# Row 113
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
# Row 122
$ws.Range("H122").Value = 1172.2858
$ws.Range("I122").Value = 1216.3077
$ws.Range("K122").Value = 3648.9231
$ws.Range("M122").Value = -1198.9231
# Row 126
$ws.Range("H126").Value = 1187.5
$ws.Range("J126").Value = 1375
$ws.Range("L126").Value = 4125
$ws.Range("N126").Value = -9065
# Row 132
$ws.Range("H132").Value = 5201.6665
$ws.Range("I132").Value = 4942
$ws.Range("K132").Value = 14826
$ws.Range("M132").Value = -12296
# Row 134
$ws.Range("H134").Value = 2200
$ws.Range("I134").Value = 2200
$ws.Range("K134").Value = 6600
$ws.Range("M134").Value = -4065

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3033.6924
$ws.Range("I5").Value = 3272
$ws.Range("J5").Value = 2927.7778
$ws.Range("K5").Value = 9816
$ws.Range("L5").Value = 8783.3334
$ws.Range("M5").Value = -9704
$ws.Range("N5").Value = -9007.3334
# Row 26
$ws.Range("H26").Value = 72.57143000000001
$ws.Range("I26").Value = 72.57143000000001
$ws.Range("K26").Value = 217.71429
$ws.Range("M26").Value = 70.28570999999999
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
# Row 32
$ws.Range("H32").Value = 1490.9231
$ws.Range("I32").Value = 1490.9231
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4472.7693
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4189.7693
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 2225
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2225
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6675
$ws.Range("N34").Value = -6843
$ws.Range("M34").ClearContents()
# Row 60
$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -29749
$ws.Range("N60").ClearContents()
# Row 103
$ws.Range("H103").Value = 4518.5
$ws.Range("J103").Value = 4337.1665
$ws.Range("L103").Value = 13011.4995
$ws.Range("N103").Value = -14769.4995
# Row 107
$ws.Range("H107").Value = 512.6
$ws.Range("J107").Value = 591.6667
$ws.Range("L107").Value = 1775.0001
$ws.Range("N107").Value = -5615.0001
# Row 135
$ws.Range("H135").Value = 3033.6924
$ws.Range("I135").Value = 3272
$ws.Range("J135").Value = 2927.7778
$ws.Range("K135").Value = 29448
$ws.Range("L135").Value = 26350.0002
$ws.Range("M135").Value = -26913
$ws.Range("N135").Value = -31420.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
# Row 132
$ws.Range("H132").Value = 5438.625
$ws.Range("I132").Value = 6070.1665
$ws.Range("J132").Value = 3544
$ws.Range("K132").Value = 18210.4995
$ws.Range("L132").Value = 10632
$ws.Range("M132").Value = -15680.4995
$ws.Range("N132").Value = -15692

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 935.6667
$ws.Range("I9").Value = 808
$ws.Range("J9").Value = 999.5
$ws.Range("K9").Value = 808
$ws.Range("L9").Value = 999.5
$ws.Range("M9").Value = -584
$ws.Range("N9").Value = -1447.5
# Row 16
$ws.Range("H16").Value = 1652.8572
$ws.Range("I16").Value = 1652.8572
$ws.Range("K16").Value = 1652.8572
$ws.Range("M16").Value = -1482.8572
# Row 46
$ws.Range("H46").Value = 798.3333
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 127
$ws.Range("H127").Value = 79999
$ws.Range("J127").Value = 79999
$ws.Range("L127").Value = 79999
$ws.Range("N127").Value = -89919

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 31002
$ws.Range("I14").Value = 31002
$ws.Range("K14").Value = 31002
$ws.Range("M14").Value = -30834
# Row 81
$ws.Range("H81").Value = 2084
$ws.Range("I81").Value = 1730.0834
$ws.Range("K81").Value = 3460.1668
$ws.Range("M81").Value = -2399.1668
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
# Row 84
$ws.Range("H84").Value = 2084
$ws.Range("I84").Value = 1730.0834
$ws.Range("K84").Value = 17300.834
$ws.Range("M84").Value = -11996.834
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
# Row 132
$ws.Range("H132").Value = 4102
$ws.Range("I132").Value = 4102
$ws.Range("K132").Value = 12306
$ws.Range("M132").Value = -9776
# Row 136
$ws.Range("H136").Value = 2456.6667
$ws.Range("I136").Value = 2483.5293
$ws.Range("K136").Value = 7450.5879
$ws.Range("M136").Value = -4900.5879
